$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B26").Value = "Are there any new FASB rulings to get up to speed on?"
$ws.Range("B28").Value = "Did any bank statements (not online ones) arrive in mail today?"

$ws.Range("B28").Select()
